# Sort the comma-separated "Recorded By" values (column G) alphabetically
# (ASCII/ordinal order) in the "Session Analysis Results" sheet.

function OrdinalSort($arr) {
    $n = $arr.Count
    for ($i = 1; $i -lt $n; $i++) {
        $key = $arr[$i]
        $j = $i - 1
        while ($j -ge 0 -and $arr[$j].CompareTo($key) -gt 0) {
            $arr[$j + 1] = $arr[$j]
            $j = $j - 1
        }
        $arr[$j + 1] = $key
    }
    return $arr
}

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -gt 1) {
            $sorted = OrdinalSort $parts
            $newVal = [string]::Join(", ", $sorted)
            if ($newVal -ne $val) {
                $cell.Value2 = $newVal
            }
        }
    }
}
